$wb = $excel.ActiveWorkbook

# Sheet ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1205.5264
$ws.Range("I32").Value = 938.53845
$ws.Range("J32").Value = 1784
$ws.Range("K32").Value = 938.53845
$ws.Range("L32").Value = 1784
$ws.Range("M32").Value = -612.53845
$ws.Range("N32").Value = -2436

# Sheet ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 71431760
$ws.Range("I62").Value = 83336216
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 83336216
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -83335592
$ws.Range("N62").Value = -6248

# Sheet ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 71431760
$ws.Range("I65").Value = 83336216
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 416681080
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -416677960
$ws.Range("N65").Value = -31240

# Sheet ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1762.7142
$ws.Range("I112").Value = 199.5
$ws.Range("J112").Value = 2388
$ws.Range("K112").Value = 598.5
$ws.Range("L112").Value = 7164
$ws.Range("M112").Value = 509.5
$ws.Range("N112").Value = -9380

# Sheet ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1022
$ws.Range("I125").Value = 1022
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 9198
$ws.Range("L125").ClearContents()
$ws.Range("M125").Value = -6738
$ws.Range("N125").ClearContents()

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18523432
$ws.Range("I32").Value = 21741290
$ws.Range("J32").Value = 20739.25
$ws.Range("K32").Value = 21741290
$ws.Range("L32").Value = 20739.25
$ws.Range("M32").Value = -21741003
$ws.Range("N32").Value = -21313.25

# Sheet ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3641.2
$ws.Range("I61").Value = 1485.3334
$ws.Range("J61").Value = 6875
$ws.Range("K61").Value = 1485.3334
$ws.Range("L61").Value = 6875
$ws.Range("M61").Value = -1273.3334
$ws.Range("N61").Value = -7299

# Sheet ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 861
$ws.Range("I97").Value = 842.6
$ws.Range("J97").Value = 999
$ws.Range("K97").Value = 842.6
$ws.Range("L97").Value = 999
$ws.Range("M97").Value = -346.6
$ws.Range("N97").Value = -1991

# Sheet ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1090.8334
$ws.Range("I110").Value = 553.63635
$ws.Range("J110").Value = 7000
$ws.Range("K110").Value = 553.63635
$ws.Range("L110").Value = 7000
$ws.Range("M110").Value = 1491.36365
$ws.Range("N110").Value = -11090

# Sheet ARM row 113
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 52299
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 52299
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 52299
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -60977

# Sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1462.8667
$ws.Range("I122").Value = 1271.8889
$ws.Range("J122").Value = 1749.3334
$ws.Range("K122").Value = 3815.6667
$ws.Range("L122").Value = 5248.0002
$ws.Range("M122").Value = -1365.6667
$ws.Range("N122").Value = -10148.0002

# Sheet ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3641.2
$ws.Range("I136").Value = 1485.3334
$ws.Range("J136").Value = 6875
$ws.Range("K136").Value = 4456.0002
$ws.Range("L136").Value = 20625
$ws.Range("M136").Value = -1906.0002
$ws.Range("N136").Value = -25725

# Sheet BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 241.57143
$ws.Range("I22").Value = 238.5
$ws.Range("J22").Value = 260
$ws.Range("K22").Value = 238.5
$ws.Range("L22").Value = 260
$ws.Range("M22").Value = -65.5
$ws.Range("N22").Value = -606

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2878.244
$ws.Range("I31").Value = 2172.8667
$ws.Range("J31").Value = 4802
$ws.Range("K31").Value = 2172.8667
$ws.Range("L31").Value = 4802
$ws.Range("M31").Value = -1877.8667
$ws.Range("N31").Value = -5392

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2878.244
$ws.Range("I34").Value = 2172.8667
$ws.Range("J34").Value = 4802
$ws.Range("K34").Value = 2172.8667
$ws.Range("L34").Value = 4802
$ws.Range("M34").Value = -1970.8667
$ws.Range("N34").Value = -5206

# Sheet CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1586.6086
$ws.Range("I122").Value = 1549.4286
$ws.Range("J122").Value = 1644.4445
$ws.Range("K122").Value = 4648.2858
$ws.Range("L122").Value = 4933.333500000001
$ws.Range("M122").Value = -2198.2858
$ws.Range("N122").Value = -9833.333500000001

# Sheet CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 568
$ws.Range("I5").Value = 470.08334
$ws.Range("J5").Value = 959.6667
$ws.Range("K5").Value = 1410.25002
$ws.Range("L5").Value = 2879.0001
$ws.Range("M5").Value = -1298.25002
$ws.Range("N5").Value = -3103.0001

# Sheet CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1437330
$ws.Range("I113").Value = 3831797.8
$ws.Range("J113").Value = 649.4
$ws.Range("K113").Value = 11495393.4
$ws.Range("L113").Value = 1948.2
$ws.Range("M113").Value = -11493223.4
$ws.Range("N113").Value = -6288.2

# Sheet CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3702.125
$ws.Range("I117").Value = 5614
$ws.Range("J117").Value = 3429
$ws.Range("K117").Value = 16842
$ws.Range("L117").Value = 10287
$ws.Range("M117").Value = -13400
$ws.Range("N117").Value = -17171

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 761.69354
$ws.Range("I131").Value = 422.5
$ws.Range("J131").Value = 923.2143
$ws.Range("K131").Value = 1267.5
$ws.Range("L131").Value = 2769.6429
$ws.Range("M131").Value = 3772.5
$ws.Range("N131").Value = -12849.6429

# Sheet CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1264213.8
$ws.Range("J132").Value = 1685168.4
$ws.Range("L132").Value = 15166515.6
$ws.Range("N132").Value = -15171575.6

# Sheet CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 568
$ws.Range("I135").Value = 470.08334
$ws.Range("J135").Value = 959.6667
$ws.Range("K135").Value = 4230.75006
$ws.Range("L135").Value = 8637.0003
$ws.Range("M135").Value = -1695.75006
$ws.Range("N135").Value = -13707.0003

# Sheet GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 337.85715
$ws.Range("I107").Value = 344.16666
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 344.16666
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1575.83334
$ws.Range("N107").Value = -4140

# Sheet LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 38846600
$ws.Range("I22").Value = 56111412
$ws.Range("J22").Value = 779.25
$ws.Range("K22").Value = 56111412
$ws.Range("L22").Value = 779.25
$ws.Range("M22").Value = -56111117
$ws.Range("N22").Value = -1369.25

# Sheet LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 38846600
$ws.Range("I27").Value = 56111412
$ws.Range("J27").Value = 779.25
$ws.Range("K27").Value = 56111412
$ws.Range("L27").Value = 779.25
$ws.Range("M27").Value = -56111305
$ws.Range("N27").Value = -993.25

# Sheet LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 10107.923
$ws.Range("I61").Value = 11218.454
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 11218.454
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -11016.454
$ws.Range("N61").Value = -4404

# Sheet LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 10107.923
$ws.Range("I113").Value = 11218.454
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 11218.454
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -9048.454
$ws.Range("N113").Value = -8340

# Sheet LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 876
$ws.Range("I122").Value = 876
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2628
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = -178
$ws.Range("N122").ClearContents()

# Sheet WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 41667588
$ws.Range("I107").Value = 83333736
$ws.Range("J107").Value = 1440
$ws.Range("K107").Value = 250001208
$ws.Range("L107").Value = 4320
$ws.Range("M107").Value = -249999288
$ws.Range("N107").Value = -8160

# Sheet WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 603.75
$ws.Range("I113").Value = 316.66666
$ws.Range("J113").Value = 1465
$ws.Range("K113").Value = 949.9999799999999
$ws.Range("L113").Value = 4395
$ws.Range("M113").Value = 1220.00002
$ws.Range("N113").Value = -8735

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2370.1667
$ws.Range("I132").Value = 1145.7333
$ws.Range("J132").Value = 4410.8887
$ws.Range("K132").Value = 3437.199900000001
$ws.Range("L132").Value = 13232.6661
$ws.Range("M132").Value = -907.1999000000005
$ws.Range("N132").Value = -18292.6661

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1628.4318
$ws.Range("I136").Value = 1365.7878
$ws.Range("J136").Value = 2416.3635
$ws.Range("K136").Value = 4097.3634
$ws.Range("L136").Value = 7249.0905
$ws.Range("M136").Value = -1547.3634
$ws.Range("N136").Value = -12349.0905

